$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.561645865440369
$ws.Range("B1").Value = 1.745550513267517
$ws.Range("C1").Value = 2.115542650222778
$ws.Range("D1").Value = 2.30855655670166
$ws.Range("E1").Value = 1.441492676734924
